# Leave Card update for FELICIDARIO, PAMELA - 12/22/2023 10:59 AM
# - Adds a new "2024" year-header row to the leave table (table grows by one row)
# - Records a new SL(1-0-0) (1.25 day) leave entry spanning three pay periods
#   (rows that used to be 74-76), including a new entry in the previously-empty
#   row 75 with its date range end (column K)
# - Recomputes the summary BALANCE cells (E9, I9) and table row formulas
# - Updates the last saved cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# --- Insert a new row above the current row 78 -----------------------------
# This shifts every row from 78 downward by one (old row 78 -> 79, ...,
# old row 142 -> 143), automatically carrying each row's original formatting
# (including the special "final row" border style that lived on row 142,
# which now lives on row 143).
$ws.Rows.Item(78).Insert() | Out-Null

# Grow Table1 so it covers the newly inserted row as well (A8:K142 -> A8:K143)
$tbl.Resize($ws.Range("A8:K143")) | Out-Null

# The inserted row does not inherit formatting automatically, so copy it from
# the data row right above (row 77) first ...
$ws.Range("A77:K77").Copy() | Out-Null
$ws.Range("A78:K78").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ... then restyle column A of the new row like the other yearly headers
# (e.g. A17 = "2020") so it gets the bold/centered "year" look (style with
# quoted-text date format), and set its text to "2024".
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A78").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A78").Value = "'2024"

# Re-apply the calculated "EARNED " helper-column formula on the new row and
# on the new final row (both lost/altered their formula during the insert).
$ws.Range("G78").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G143").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- New SL(1-0-0) leave entries (1.25 day each) ---------------------------
# Row 74 (9/1/2023 period) already had the particulars/undertime filled in;
# only the EARNED amount was missing.
$ws.Range("C74").Value = 1.25

# Row 75 (10/1/2023 period) is a brand-new entry: particulars, earned amount,
# the "W/ Pay" undertime flag and the covered date (column K).
$ws.Range("B75").Value = "SL(1-0-0)"
$ws.Range("C75").Value = 1.25
$ws.Range("H75").Value = 1
$ws.Range("K74").Copy() | Out-Null
$ws.Range("K75").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("K75").Value = 45202

# Row 76 (11/1/2023 period) only needed the EARNED amount.
$ws.Range("C76").Value = 1.25

# --- Recalculate everything -------------------------------------------------
$excel.CalculateFull() | Out-Null

# --- Restore the last active selection -------------------------------------
$ws.Range("B68").Select() | Out-Null
